# Update imputed KNN result values (Name of Algo update)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A3").Value = -21.937
$ws.Range("E3").Value = 16.483
$ws.Range("A21").Value = -19.861
$ws.Range("A23").Value = -20.317
$ws.Range("E24").Value = 16.622
$ws.Range("A25").Value = -21.765
$ws.Range("D27").Value = -8.388999999999999
$ws.Range("D31").Value = -8.247999999999999
$ws.Range("D39").Value = -7.553
$ws.Range("D48").Value = -7.475
$ws.Range("D51").Value = -8.374000000000001
$ws.Range("D52").Value = -8.083
$ws.Range("A53").Value = -21.993
$ws.Range("D55").Value = -8.035
$ws.Range("D56").Value = -8.409000000000001
$ws.Range("A57").Value = -22.563
$ws.Range("D57").Value = -8.106999999999999
$ws.Range("E57").Value = 16.574
$ws.Range("A59").Value = -22.5
$ws.Range("E61").Value = 16.602
$ws.Range("A69").Value = -21.649
$ws.Range("E70").Value = 17.696
$ws.Range("D73").Value = -8.004000000000001
$ws.Range("A79").Value = -20.919
$ws.Range("A83").Value = -21.997
$ws.Range("E86").Value = 16.597
$ws.Range("D89").Value = -6.667
$ws.Range("D90").Value = -7.475999999999999
$ws.Range("A93").Value = -21.472
$ws.Range("E98").Value = 16.421
$ws.Range("E100").Value = 16.725
$ws.Range("E102").Value = 16.49
